# Ancient Gear and Darklord update
# Adds a new "SR03-JP" sheet (Structure Deck R: Machine Dragon Re-Volt)
# and nudges the saved selection on the three existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet after the last existing sheet (SPDS-JP)
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$newSheet   = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "SR03-JP"

# ---------------------------------------------------------------------
# 2. Fill in the card-id rows. Column A is filled top-to-bottom first so
#    the new shared strings land in the same order as the source sheet,
#    then the numeric id column, then the separator columns.
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Structure Deck R: Machine Dragon Re-Volt"
$newSheet.Range("A2").Value = "Ancient Gear Gadget"
$newSheet.Range("A3").Value = "Ancient Gear Reactor Dragon"
$newSheet.Range("A4").Value = "Ancient Gear Catapult"

$newSheet.Range("B2").Value = 100303000
$newSheet.Range("B3").Value = 100303001
$newSheet.Range("B4").Value = 100303021

$newSheet.Range("C2").Value = ":"
$newSheet.Range("C3").Value = ":"
$newSheet.Range("C4").Value = ":"

$newSheet.Range("E2").Value = ";"
$newSheet.Range("E3").Value = ";"
$newSheet.Range("E4").Value = ";"

# ---------------------------------------------------------------------
# 3. B1 holds the release date. Clone the date-header formatting
#    (medium border + numFmtId 15) from an existing sheet's header cell,
#    then drop the bold so it matches the other header cells' look, and
#    finally set the date serial for 2016-09-24.
# ---------------------------------------------------------------------
$dateTemplate = $wb.Worksheets.Item("VP16-JP").Range("B1")
$dateTemplate.Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$newSheet.Range("B1").Font.Bold = $false
$newSheet.Range("B1").Value = 42637

# ---------------------------------------------------------------------
# 4. Update the remembered selection on each sheet.
# ---------------------------------------------------------------------
$null = $newSheet.Range("B8").Select()

$null = $wb.Worksheets.Item("VP16-JP").Range("F8").Select()
$null = $wb.Worksheets.Item("TDIL-EN").Range("A16").Select()

$spds = $wb.Worksheets.Item("SPDS-JP")
$null = $spds.Activate()
$null = $spds.Range("B34").Select()
